# Performance benchmark results refresh for collect.xlsx
# Updates benchmark timing figures on the "collect_with_push" (sheet1) and
# "collect_with_extend" (sheet2) sheets, turns the "rayon" baseline rows on
# sheet1 into cross-sheet references to sheet2 (rayon's numbers are shared),
# and moves the active tab / selection to match the author's final state.

$wb = $excel.ActiveWorkbook

$wsPush   = $wb.Worksheets.Item("collect_with_push")
$wsExtend = $wb.Worksheets.Item("collect_with_extend")

# ---------------------------------------------------------------------------
# Sheet "collect_with_push" (sheet1): updated raw benchmark values
# ---------------------------------------------------------------------------

# Row 2 ("rayon", 16384 elements): now mirrors collect_with_extend row 2
$wsPush.Range("E2").Formula = "=collect_with_extend!E2"
$wsPush.Range("G2").Formula = "=collect_with_extend!G2"

# Row 3 (SplitVec Doubling push, 16384 elements)
$wsPush.Range("G3").Value = 18.974

# Row 4 (SplitVec Linear push, 16384 elements)
$wsPush.Range("E4").Value = 4.4450000000000003
$wsPush.Range("G4").Value = 14.282999999999999

# Row 5 (FixedVec push, 16384 elements)
$wsPush.Range("E5").Value = 4.1616
$wsPush.Range("G5").Value = 12.961

# Row 10 ("rayon", 65536 elements): now mirrors collect_with_extend row 10
$wsPush.Range("E10").Formula = "=collect_with_extend!E10"
$wsPush.Range("G10").Formula = "=collect_with_extend!G10"

# Row 11 (SplitVec Doubling push, 65536 elements)
$wsPush.Range("G11").Value = 82.436999999999998

# Row 12 (SplitVec Linear push, 65536 elements)
$wsPush.Range("E12").Value = 20.071000000000002
$wsPush.Range("G12").Value = 65.004000000000005

# Row 13 (FixedVec push, 65536 elements)
$wsPush.Range("E13").Value = 19.440999999999999
$wsPush.Range("G13").Value = 48.92

# ---------------------------------------------------------------------------
# Sheet "collect_with_extend" (sheet2): updated raw benchmark values
# ---------------------------------------------------------------------------

# Row 3 (SplitVec Doubling extend, 16384 elements)
$wsExtend.Range("E3").Value = 0.59357000000000004
$wsExtend.Range("G3").Value = 16.343

# Row 4 (SplitVec Linear extend, 16384 elements)
$wsExtend.Range("E4").Value = 0.35904999999999998
$wsExtend.Range("G4").Value = 9.4992000000000001

# Row 5 (FixedVec extend, 16384 elements)
$wsExtend.Range("E5").Value = 0.34448000000000001
$wsExtend.Range("G5").Value = 6.4676

# Row 11 (SplitVec Doubling extend, 65536 elements)
$wsExtend.Range("E11").Value = 2.1206999999999998
$wsExtend.Range("G11").Value = 64.709000000000003

# Row 12 (SplitVec Linear extend, 65536 elements)
$wsExtend.Range("E12").Value = 1.0075000000000001
$wsExtend.Range("G12").Value = 41.076999999999998

# Row 13 (FixedVec extend, 65536 elements)
$wsExtend.Range("E13").Value = 1.0028999999999999
$wsExtend.Range("G13").Value = 26.298999999999999

# ---------------------------------------------------------------------------
# Selections / active tab: collect_with_push loses focus, selects C44;
# collect_with_extend becomes the active tab, selecting C45:C46.
# ---------------------------------------------------------------------------

$wsPush.Range("C44").Select()

$wsExtend.Activate()
$wsExtend.Range("C45:C46").Select()
